$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9 (this shifts the old rows 9-14 down to 10-15),
# to make room for the new "12:20" period that was missing from the schedule.
$ws.Rows.Item(9).Insert()

# Two extra periods (17:30 and 18:20) are appended at the end of the day (rows 16-17)
# so the schedule covers a full 6 hour afternoon turn.
$ws.Rows.Item(16).Insert()
$ws.Rows.Item(17).Insert()

# --- Column A: time-of-day labels for every row ---
$times = @("ELT-2A","7:00","7:50","8:40","9:30","9:50","10:40","11:30","12:20","13:00","13:50","14:40","15:30","15:50","16:40","17:30","18:20")
for ($i = 0; $i -lt $times.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $times[$i]
}

# --- Row 1 header (days of week), unchanged ---
$ws.Range("B1").Value = "segunda"
$ws.Range("C1").Value = "terça"
$ws.Range("D1").Value = "quarta"
$ws.Range("E1").Value = "quinta"
$ws.Range("F1").Value = "sexta"

# --- Full schedule grid, rows 2-16 (columns B:F) ---
$grid = @(
    @("-","-","-","-","-"),                                                                                  # row 2  - 7:00
    @("-","-","-","-","Lucas Ferreira-Sistemas digitais"),                                                    # row 3  - 7:50
    @("-","-","Josivaldo Ferreira-Circuitos Elétricos 2","-","Lucas Ferreira-Sistemas digitais"),             # row 4  - 8:40
    @("Intervalo","Intervalo","Intervalo","Intervalo","Intervalo"),                                            # row 5  - 9:30
    @("Josivaldo Ferreira-Circuitos Elétricos 2","-","Josivaldo Ferreira-Circuitos Elétricos 2","-","-"),     # row 6  - 9:50
    @("-","Andre Lucca-Acionamentos","Andre Lucca-Acionamentos","-","-"),                                      # row 7  - 10:40
    @("-","-","-","-","-"),                                                                                   # row 8  - 11:30
    @("Almoço","Almoço","Almoço","Almoço","Almoço"),                                                          # row 9  - 12:20
    @("-","-","-","-","-"),                                                                                   # row 10 - 13:00
    @("-","-","-","-","-"),                                                                                   # row 11 - 13:50
    @("-","-","-","-","-"),                                                                                   # row 12 - 14:40
    @("Intervalo","Intervalo","Intervalo","Intervalo","Intervalo"),                                            # row 13 - 15:30
    @("-","-","-","-","-"),                                                                                   # row 14 - 15:50
    @("-","-","-","-","-"),                                                                                   # row 15 - 16:40
    @("-","-","-","-","-")                                                                                    # row 16 - 17:30
)

$colLetters = @("B","C","D","E","F")
for ($r = 0; $r -lt $grid.Length; $r++) {
    $rowValues = $grid[$r]
    $rowNumber = $r + 2
    for ($c = 0; $c -lt $colLetters.Length; $c++) {
        $ws.Cells.Item($rowNumber, $c + 2).Value = $rowValues[$c]
    }
}

# --- Row 17 (18:20) is left blank in columns B:F ---
$ws.Range("B17").Value = ""
$ws.Range("C17").Value = ""
$ws.Range("D17").Value = ""
$ws.Range("E17").Value = ""
$ws.Range("F17").Value = ""
